# Updates cryptos list figures (prices / 1h volume change %) to the
# latest scrape, matching the Sat Sep 23 07:30:44 UTC 2023 GitHub Actions
# run.
#
# The Price column (D) stores its figures as literal text (e.g.
# "26.647.61", "1.598.16" -- these use '.' as a thousands separator, not
# a decimal point, so they are not valid numbers). Whenever a new Price
# value happens to look like a genuine number (e.g. "211.26"), assigning
# it straight to .Value would make Excel reinterpret it as a numeric
# cell, which would both change the cell's type and round/reformat the
# value. To keep it as text we enter it the same way a user would in the
# Excel UI -- with a leading apostrophe -- and then reapply the "Normal"
# style so the cell doesn't pick up a lingering quote-prefix indicator.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.651.76"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.598.12"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'211.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("E6").Value = "  +0.76%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").Value = "'19.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("D11").Value = "'0.0840"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "1.822.41"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "1.603.24"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").Value = "'4.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").Value = "'0.523"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").Value = "'64.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "26.644.31"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "0.0₃0733"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "'207.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.53%  "
$ws.Range("E21").Value = "  +5.73%  "
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("E23").Value = "  +1.14%  "
$ws.Range("D24").Value = "'8.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'145.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").Value = "'15.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.275.61"
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.622"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.25%  "
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("D37").Value = "'1.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("D38").Value = "'0.0170"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("E39").Value = "  +20.64%  "
$ws.Range("D40").Value = "'0.834"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("E41").Value = "  +2.80%  "
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").Value = "'0.785"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").Value = "'63.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("D45").Value = "1.735.19"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "'90.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("D47").Value = "'1.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.26%  "
$ws.Range("E48").Value = "  +3.42%  "
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").Value = "'7.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.73%  "
